$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Kevin Blenman's password (row 7, column B) to "tigers"
$ws.Range("B7").Value = "tigers"

# Update the active selection to D12
$ws.Range("D12").Select()
